$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture formats we will need later (before the source cells are cleared
#    or shifted by the row insertion below).
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("X7").PasteSpecial(-4122) | Out-Null

$ws.Range("A7").Copy() | Out-Null
$ws.Range("X8").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("X2").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a new row for the "Node Status" UUID entry (pushes rows 11-23
#    down to 12-24, inheriting formatting from the row above - row 9).
# ---------------------------------------------------------------------------
$ws.Rows("10").Insert()

# ---------------------------------------------------------------------------
# 3. Remove the old UUID values that used to live in column B, now that
#    their formats have been copied to column X.
# ---------------------------------------------------------------------------
$ws.Range("B1").ClearContents()
$ws.Range("B7").ClearContents()

# ---------------------------------------------------------------------------
# 4. Populate the new "UUID" column (X) with header + values.
# ---------------------------------------------------------------------------
$ws.Range("X2").Value = "UUID"
$ws.Range("X7").Value = "4fafc201-1fb5-459e-8fcc-c5c9c331914b"
$ws.Range("X8").Value = "beb5483e-36e1-4688-b7f5-ea07361b26a8"

# ---------------------------------------------------------------------------
# 5. Fill in the newly inserted row 10 ("Node Status" characteristic).
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "Node Status"
$ws.Range("C10").Value = "Read | Notify"

# ---------------------------------------------------------------------------
# 6. Column X width.
# ---------------------------------------------------------------------------
$ws.Range("X1").ColumnWidth = 47.71

# ---------------------------------------------------------------------------
# 7. Selection as saved in the file.
# ---------------------------------------------------------------------------
$ws.Range("D10").Select() | Out-Null
